$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.017517448585473
$ws.Range("C2").Value = 0.4223807153183843
$ws.Range("D2").Value = 0.2446916143694868
$ws.Range("E2").Value = 1.60357238366937
$ws.Range("F2").Value = -0.1128288494190266
$ws.Range("G2").Value = 0.564231488800237

$ws.Range("B3").Value = 1.676451378490478
$ws.Range("C3").Value = 0.6391979074220058
$ws.Range("D3").Value = 0.2445720448915195
$ws.Range("E3").Value = 2.490873017364748
$ws.Range("F3").Value = -0.1127528199373788
$ws.Range("G3").Value = 0.5639702106451753

$ws.Range("B4").Value = 0.8448534768133353
$ws.Range("C4").Value = 0.3599118058759405
$ws.Range("D4").Value = 0.2127515357305201
$ws.Range("E4").Value = 1.550690596630148
$ws.Range("F4").Value = -0.09281505800686342
$ws.Range("G4").Value = 0.4989226723653052
